# ---------------------------------------------------------------------------
# Edit: (1) re-style the table on slide 16 with a different built-in table
#       style id, and (2) swap the presentation's applied colour theme from
#       "Integral" to the default "Office Theme" palette (the deck's table
#       styles/colours live in ppt/theme/theme2.xml, the theme actually
#       bound to the slide master that the slides use).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style change (slide 16, the graphicFrame holding the table) --
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{49F01519-4898-4B14-80E6-56CDCBF296A0}")

# --- 2) Swap the theme colour scheme to the stock "Office Theme" colours --
# dk1 (000000) and lt1 (FFFFFF) are identical between the two themes, so
# only the remaining ten slots need to change.
$officeThemeColors = @{
    3  = 0x44546A   # dk2
    4  = 0xE7E6E6   # lt2
    5  = 0x5B9BD5   # accent1
    6  = 0xED7D31   # accent2
    7  = 0xA5A5A5   # accent3
    8  = 0xFFC000   # accent4
    9  = 0x4472C4   # accent5
    10 = 0x70AD47   # accent6
    11 = 0x0563C1   # hyperlink
    12 = 0x954F72   # followed hyperlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($slotIndex in $officeThemeColors.Keys) {
    $hex = $officeThemeColors[$slotIndex]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint's .RGB long is packed little-endian (0x00BBGGRR).
    $bgr = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($slotIndex).RGB = $bgr
}
